{"js": "// Remove the hyperparameters/momentum write-up text from its paragraph\n// (keeping the empty paragraph + its _GoBack bookmark), and remove the\n// now-redundant blank indented paragraph that followed it.\n\nconst body = context.document.body;\n\n// Locate the paragraph that contains the draft text by searching for a\n// distinctive, stable substring rather than relying on a hard-coded index.\nconst searchResults = body.search(\"The hyperparameters we used for\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const hit = searchResults.items[0];\n  const hitParagraphs = hit.paragraphs;\n  hitParagraphs.load(\"items\");\n  await context.sync();\n\n  const targetParagraph = hitParagraphs.items[0];\n\n  // Grab the next paragraph (the blank indented one) before we mutate\n  // the target paragraph, so we can remove it afterwards.\n  const nextParagraph = targetParagraph.getNextOrNullObject();\n  nextParagraph.load(\"isNullObject,text\");\n\n  // Clear just the text of the target paragraph's range; this removes the\n  // runs but preserves the paragraph mark and the bookmarkStart/bookmarkEnd\n  // pair (unlike Paragraph.clear(), which would also drop the bookmark).\n  const targetRange = targetParagraph.getRange();\n  targetRange.insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Remove the following blank paragraph entirely, matching the diff.\n  if (!nextParagraph.isNullObject && nextParagraph.text === \"\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the hyperparameters/momentum write-up text from its paragraph\n# (keeping the empty paragraph + its _GoBack bookmark), and remove the\n# now-redundant blank indented paragraph that followed it.\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"The hyperparameters we used for\")\n\nif ($found) {\n    # Expand to the whole paragraph containing the match.\n    [void]$range.Expand(4)  # wdParagraph\n\n    # Grab the following paragraph (the blank indented one) before editing,\n    # so we still have a handle on it afterwards.\n    $nextRange = $range.Next(4, 1)  # wdParagraph, 1 unit forward\n    if ($nextRange -ne $null) {\n        [void]$nextRange.Expand(4)\n    }\n\n    # Drop the trailing paragraph mark from the selection so clearing the\n    # text doesn't touch the paragraph mark / bookmarkStart-bookmarkEnd\n    # pair that must remain.\n    [void]$range.MoveEnd(1, -1)  # wdCharacter\n    $range.Text = \"\"\n\n    # Delete the following blank paragraph entirely, matching the diff.\n    # (A paragraph range's .Text always ends with its paragraph mark, so an\n    # \"empty\" paragraph reads back as just \"`r\" rather than \"\".)\n    if ($nextRange -ne $null -and $nextRange.Text.Trim() -eq \"\") {\n        $nextRange.Delete()\n    }\n}\n"}
